# Apply the commit's changes to the workbook.
#
# Summary of the edit (derived from the OOXML diff):
#   - Sheet 1 ("展览") and Sheet 4 ("全部类型") each had their row 2
#     (the "2024-04-21 苏州·梦幻岛..." entry) removed, shifting every
#     subsequent row up by one.
#   - Column A holds a plain static sequence number (row-1) that is
#     NOT part of the shifted content in the diff - it keeps its
#     original values (1, 2, 3, ...) after the row removal, so it
#     must be restored after the row delete (which otherwise shifts
#     it along with everything else).
#   - After the shift, a handful of rows also got their "想去人数"
#     (column F, want-to-go count) value bumped to a newer number,
#     reflecting re-scraped data.
#   - Sheets 2 ("演出") and 3 ("本地生活") are unaffected.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Delete row 2 on sheet 1 and sheet 4 (shifts subsequent rows up)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Delete()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Delete()

# ------------------------------------------------------------------
# 1b) Restore column A's static sequence numbers (row - 1), since
#     the row delete above also shifted those values up, which is
#     not what the source diff shows.
# ------------------------------------------------------------------
$used1 = $ws1.UsedRange.Rows.Count
for ($r = 1; $r -le $used1; $r++) {
    $ws1.Cells.Item($r, 1).Value2 = $r - 1
}

$used4 = $ws4.UsedRange.Rows.Count
for ($r = 1; $r -le $used4; $r++) {
    $ws4.Cells.Item($r, 1).Value2 = $r - 1
}

# ------------------------------------------------------------------
# 2) Update the "want-to-go count" (column F) on specific rows that
#    changed value beyond the pure row-shift, for each sheet.
#    Each tuple is (rowNumber, expectedOldValue, newValue); the old
#    value is only used as a sanity check.
# ------------------------------------------------------------------
$sheet1Changes = @(
    ,@(2, 14781, 14794)
    ,@(3, 18208, 18238)
    ,@(5, 99, 100)
    ,@(13, 49, 50)
    ,@(14, 76, 81)
    ,@(15, 188, 189)
    ,@(17, 1376, 1379)
    ,@(20, 77, 79)
    ,@(21, 222, 223)
    ,@(22, 7534, 7553)
    ,@(24, 13, 14)
    ,@(25, 48, 49)
    ,@(26, 1198, 1200)
    ,@(28, 5911, 5914)
    ,@(29, 86, 87)
    ,@(30, 51, 54)
    ,@(33, 251, 253)
    ,@(34, 5223, 5236)
)

$sheet4Changes = @(
    ,@(2, 14781, 14794)
    ,@(3, 18208, 18238)
    ,@(5, 99, 100)
    ,@(13, 49, 50)
    ,@(14, 76, 81)
    ,@(15, 188, 189)
    ,@(17, 1376, 1380)
    ,@(21, 77, 79)
    ,@(22, 222, 223)
    ,@(23, 7534, 7553)
    ,@(25, 13, 14)
    ,@(26, 48, 49)
    ,@(27, 1198, 1200)
    ,@(30, 5911, 5914)
    ,@(31, 86, 87)
    ,@(32, 51, 54)
    ,@(35, 251, 253)
    ,@(36, 5223, 5236)
)

foreach ($ch in $sheet1Changes) {
    $ws1.Cells.Item($ch[0], 6).Value2 = $ch[2]
}

foreach ($ch in $sheet4Changes) {
    $ws4.Cells.Item($ch[0], 6).Value2 = $ch[2]
}
